$wb = $excel.ActiveWorkbook

# --- Generate Report for Handback ---
# The localization status report is being regenerated: the handback for
# zh-cn finished (so its "Ready for handoff" status becomes "Handed back:
# in sync with en-US" everywhere it appears), the handback timestamps for
# both the zh-cn and de-de rows move forward, and the de-de row's stale
# "handback not latest" error clears now that it is back in sync.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview summary row (zh-cn/de-de columns)
#    and the Status column on each language sheet, so replace it workbook-wide.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", 1)
}

# 2. Latest Handback DateTime refreshed for both languages.
$wsZh.Range("K2").Value = "2016-09-03 02:53:48"
$wsDe.Range("K2").Value = "2016-09-03 02:53:55"

# 3. Error Detail cleared now that the handback files are in sync.
$wsZh.Range("P2").Value = ""
$wsDe.Range("P2").Value = ""

# 4. Column widths widen to fit the longer status text (Status columns) and
#    shrink the now-empty Error Detail column back down.
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(16).ColumnWidth = 12.85

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(16).ColumnWidth = 12.85
